$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Despezas")

# Update row 2 values
$ws.Range("A2").Value = "R$ 900"
$ws.Range("B2").Value = "Teste"

# Remove rows 3 through 5 entirely (they held Cavalo vendido, Mercado, Pix tio paulo)
$ws.Rows("3:5").Delete()
